$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("départements")
$ws.Range("C6").Value = 19.23076923076923
$ws.Range("D6").Value = 5
$ws.Range("C13").Value = 9.67741935483871
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 62
$ws.Range("C19").Value = 18.42105263157895
$ws.Range("D19").Value = 7
$ws.Range("C71").Value = 42.85714285714285
$ws.Range("D71").Value = 42
$ws.Range("C72").Value = 33.33333333333333
$ws.Range("D72").Value = 12
$ws.Range("C95").Value = 69.04761904761905
$ws.Range("D95").Value = 29
$ws.Range("C103").Value = 6.451612903225806
$ws.Range("D103").Value = 2
$ws.Range("C104").Value = 6.299212598425196
$ws.Range("D104").Value = 8
$ws.Range("C121").Value = 18.0327868852459
$ws.Range("D121").Value = 11
$ws.Range("C153").Value = 9.58904109589041
$ws.Range("D153").Value = 7
$ws.Range("C158").Value = 11.32075471698113
$ws.Range("D158").Value = 24
$ws.Range("C161").Value = 10.6145251396648
$ws.Range("D161").Value = 19
$ws.Range("C166").Value = 9.67741935483871
$ws.Range("D166").Value = 12
$ws.Range("C168").Value = 9.770114942528735
$ws.Range("D168").Value = 17
$ws.Range("C169").Value = 11.9047619047619
$ws.Range("D169").Value = 5
$ws.Range("C171").Value = 12.94117647058824
$ws.Range("D171").Value = 11
$ws.Range("C175").Value = 6.557377049180328
$ws.Range("D175").Value = 8
$ws.Range("C177").Value = 11.02941176470588
$ws.Range("D177").Value = 15
$ws.Range("C182").Value = 4.347826086956522
$ws.Range("D182").Value = 7
$ws.Range("C183").Value = 5.376344086021505
$ws.Range("D183").Value = 5
$ws.Range("C187").Value = 10
$ws.Range("D187").Value = 7
$ws.Range("C190").Value = 11.02362204724409
$ws.Range("D190").Value = 14
$ws.Range("C191").Value = 9.090909090909092
$ws.Range("D191").Value = 7
$ws.Range("C200").Value = 3.225806451612903
$ws.Range("D200").Value = 1
$ws.Range("C201").Value = 7.801418439716312
$ws.Range("D201").Value = 11
$ws.Range("C213").Value = 7.142857142857142
$ws.Range("D213").Value = 4
$ws.Range("C227").Value = 5.263157894736842
$ws.Range("D227").Value = 9
$ws.Range("C255").Value = 6.302521008403361
$ws.Range("D255").Value = 15
$ws.Range("C263").Value = 8.02919708029197
$ws.Range("D263").Value = 11
$ws.Range("C264").Value = 1.886792452830189
$ws.Range("E264").Value = 106
$ws.Range("C265").Value = 10.05917159763314
$ws.Range("D265").Value = 17
$ws.Range("C266").Value = 19.56521739130435
$ws.Range("D266").Value = 9
$ws.Range("C273").Value = 6.395348837209303
$ws.Range("D273").Value = 11
$ws.Range("C274").Value = 8.51063829787234
$ws.Range("D274").Value = 12
$ws.Range("C279").Value = 4.907975460122699
$ws.Range("D279").Value = 8
$ws.Range("C280").Value = 6.896551724137931
$ws.Range("D280").Value = 6
$ws.Range("C288").Value = 10.12658227848101
$ws.Range("D288").Value = 8
$ws.Range("C297").Value = 4.651162790697675
$ws.Range("D297").Value = 2
$ws.Range("C310").Value = 3.278688524590164
$ws.Range("D310").Value = 2
$ws.Range("C331").Value = 1.621621621621622
$ws.Range("D331").Value = 3
$ws.Range("C332").Value = 2.597402597402597
$ws.Range("D332").Value = 2
$ws.Range("C342").Value = 0.9345794392523363
$ws.Range("D342").Value = 1
$ws.Range("C353").Value = 1.724137931034483
$ws.Range("D353").Value = 2
$ws.Range("C355").Value = 0.4347826086956522
$ws.Range("E355").Value = 230
$ws.Range("C362").Value = 3.535353535353535
$ws.Range("D362").Value = 7
$ws.Range("C363").Value = 6.779661016949152
$ws.Range("D363").Value = 4
$ws.Range("C365").Value = 2.97029702970297
$ws.Range("D365").Value = 3
$ws.Range("C373").Value = 2.678571428571428
$ws.Range("D373").Value = 3
$ws.Range("C377").Value = 0.9433962264150944
$ws.Range("D377").Value = 1
$ws.Range("C384").Value = 7.2992700729927
$ws.Range("D384").Value = 10
$ws.Range("C388").Value = 7.964601769911504
$ws.Range("D388").Value = 9
$ws.Range("C491").Value = 6.666666666666667
$ws.Range("D491").Value = 1
$ws.Range("C494").Value = 25
$ws.Range("D494").Value = 5
$ws.Range("C496").Value = 25.92592592592592
$ws.Range("D496").Value = 7
$ws.Range("C498").Value = 48.38709677419355
$ws.Range("D498").Value = 15
$ws.Range("E498").Value = 31
$ws.Range("C499").Value = 16.40625
$ws.Range("E499").Value = 128
$ws.Range("C517").Value = 27.94117647058824
$ws.Range("D517").Value = 19
$ws.Range("C554").Value = 25
$ws.Range("D554").Value = 13
$ws.Range("C578").Value = 32.35294117647059
$ws.Range("D578").Value = 22
$ws.Range("C588").Value = 4.651162790697675
$ws.Range("D588").Value = 2
$ws.Range("C589").Value = 4.794520547945205
$ws.Range("D589").Value = 7
$ws.Range("C601").Value = 3.278688524590164
$ws.Range("D601").Value = 2
$ws.Range("C615").Value = 4.020100502512562
$ws.Range("D615").Value = 8
$ws.Range("C643").Value = 2.811244979919679
$ws.Range("D643").Value = 7
$ws.Range("C653").Value = 5.05050505050505
$ws.Range("D653").Value = 10
$ws.Range("C654").Value = 10.52631578947368
$ws.Range("D654").Value = 6
$ws.Range("C662").Value = 5.405405405405405
$ws.Range("D662").Value = 8
$ws.Range("C667").Value = 1.630434782608696
$ws.Range("D667").Value = 3
$ws.Range("C668").Value = 0.9433962264150944
$ws.Range("D668").Value = 1
$ws.Range("C685").Value = 4.651162790697675
$ws.Range("D685").Value = 2
$ws.Range("C686").Value = 4.794520547945205
$ws.Range("D686").Value = 7
$ws.Range("C698").Value = 3.278688524590164
$ws.Range("D698").Value = 2
$ws.Range("C719").Value = 2.162162162162162
$ws.Range("D719").Value = 4
$ws.Range("C720").Value = 2.597402597402597
$ws.Range("D720").Value = 2
$ws.Range("C730").Value = 0.9345794392523363
$ws.Range("D730").Value = 1
$ws.Range("C740").Value = 3.2
$ws.Range("D740").Value = 8
$ws.Range("C741").Value = 3.389830508474576
$ws.Range("D741").Value = 4
$ws.Range("C750").Value = 5.05050505050505
$ws.Range("D750").Value = 10
$ws.Range("C751").Value = 10.16949152542373
$ws.Range("D751").Value = 6
$ws.Range("C753").Value = 2.97029702970297
$ws.Range("D753").Value = 3
$ws.Range("C759").Value = 5.405405405405405
$ws.Range("D759").Value = 8
$ws.Range("C761").Value = 3.571428571428571
$ws.Range("D761").Value = 4
$ws.Range("C764").Value = 2.717391304347826
$ws.Range("D764").Value = 5
$ws.Range("C765").Value = 0.9433962264150944
$ws.Range("D765").Value = 1
$ws.Range("C772").Value = 7.971014492753622
$ws.Range("D772").Value = 11
$ws.Range("C776").Value = 8.849557522123893
$ws.Range("D776").Value = 10
$ws.Range("C782").Value = 4.651162790697675
$ws.Range("D782").Value = 2
$ws.Range("C795").Value = 3.278688524590164
$ws.Range("D795").Value = 2
$ws.Range("C809").Value = 3.015075376884422
$ws.Range("D809").Value = 6
$ws.Range("C847").Value = 3.535353535353535
$ws.Range("D847").Value = 7
$ws.Range("C848").Value = 6.779661016949152
$ws.Range("D848").Value = 4
$ws.Range("C862").Value = 0.9433962264150944
$ws.Range("D862").Value = 1

$ws = $wb.Worksheets.Item("régions")
$ws.Range("D2").Value = 10.85
$ws.Range("E2").Value = 89
$ws.Range("D3").Value = 32.35
$ws.Range("E3").Value = 132
$ws.Range("D5").Value = 5.45
$ws.Range("E5").Value = 48
$ws.Range("D6").Value = 56.02
$ws.Range("E6").Value = 228
$ws.Range("D7").Value = 11.31
$ws.Range("E7").Value = 97
$ws.Range("D8").Value = 8.130000000000001
$ws.Range("E8").Value = 72
$ws.Range("D10").Value = 7.01
$ws.Range("E10").Value = 62
$ws.Range("D14").Value = 1.12
$ws.Range("E14").Value = 5
$ws.Range("D15").Value = 28.46
$ws.Range("E15").Value = 72
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 21
$ws.Range("D17").Value = 1.78
$ws.Range("E17").Value = 8
$ws.Range("D18").Value = 1.11
$ws.Range("E18").Value = 5
$ws.Range("D19").Value = 1.79
$ws.Range("E19").Value = 8
$ws.Range("D20").Value = 9.73
$ws.Range("E20").Value = 46
$ws.Range("D23").Value = 1.95
$ws.Range("E23").Value = 13
$ws.Range("D24").Value = 24.54
$ws.Range("E24").Value = 107
$ws.Range("D25").Value = 8.880000000000001
$ws.Range("E25").Value = 50
$ws.Range("D26").Value = 2.56
$ws.Range("E26").Value = 17
$ws.Range("D27").Value = 1.8
$ws.Range("E27").Value = 12
$ws.Range("D28").Value = 2.42
$ws.Range("E28").Value = 16
$ws.Range("D29").Value = 7.65
$ws.Range("E29").Value = 31
$ws.Range("D38").Value = 9.470000000000001
$ws.Range("E38").Value = 61
$ws.Range("F41").Value = 797
$ws.Range("D43").Value = 4.78
$ws.Range("E43").Value = 35
$ws.Range("D44").Value = 2.5
$ws.Range("E44").Value = 20
$ws.Range("D46").Value = 2.14
$ws.Range("E46").Value = 17
$ws.Range("D47").Value = 6.51
$ws.Range("E47").Value = 46
$ws.Range("D48").Value = 28.24
$ws.Range("E48").Value = 98
$ws.Range("D52").Value = 5.53
$ws.Range("E52").Value = 43
$ws.Range("F52").Value = 778
$ws.Range("D56").Value = 9.09
$ws.Range("E56").Value = 40
$ws.Range("D59").Value = 1.26
$ws.Range("E59").Value = 7
$ws.Range("D62").Value = 1.8
$ws.Range("E62").Value = 10
$ws.Range("D74").Value = 6.13
$ws.Range("E74").Value = 53
$ws.Range("D84").Value = 20.04
$ws.Range("F84").Value = 539
$ws.Range("D87").Value = 25.3
$ws.Range("E87").Value = 150
$ws.Range("F87").Value = 593
$ws.Range("D88").Value = 5.97
$ws.Range("E88").Value = 54
$ws.Range("D90").Value = 1.19
$ws.Range("E90").Value = 13
$ws.Range("D91").Value = 1.83
$ws.Range("E91").Value = 20
$ws.Range("D92").Value = 9.67
$ws.Range("E92").Value = 97
$ws.Range("D95").Value = 1.66
$ws.Range("E95").Value = 22
$ws.Range("D96").Value = 25.89
$ws.Range("E96").Value = 196
$ws.Range("D97").Value = 6.55
$ws.Range("E97").Value = 72
$ws.Range("D98").Value = 2.18
$ws.Range("E98").Value = 29
$ws.Range("D99").Value = 1.28
$ws.Range("E99").Value = 17
$ws.Range("D100").Value = 1.81
$ws.Range("E100").Value = 24
$ws.Range("D101").Value = 4.76
$ws.Range("E101").Value = 33
$ws.Range("D102").Value = 15.41
$ws.Range("E102").Value = 57
$ws.Range("F102").Value = 370
$ws.Range("D104").Value = 1.62
$ws.Range("E104").Value = 13
$ws.Range("D105").Value = 34.24
$ws.Range("E105").Value = 113
$ws.Range("D106").Value = 6.85
$ws.Range("E106").Value = 49
$ws.Range("D107").Value = 2.85
$ws.Range("E107").Value = 23
$ws.Range("D108").Value = 1.36
$ws.Range("E108").Value = 11
$ws.Range("D109").Value = 2.61
$ws.Range("E109").Value = 21

$ws = $wb.Worksheets.Item("national")
$ws.Range("B2").Value = 8.109999999999999
$ws.Range("C2").Value = 624
$ws.Range("B3").Value = 23.6
$ws.Range("C3").Value = 978
$ws.Range("D3").Value = 4144
$ws.Range("B5").Value = 1.76
$ws.Range("C5").Value = 172
$ws.Range("D5").Value = 9791
$ws.Range("B6").Value = 28.64
$ws.Range("C6").Value = 1521
$ws.Range("D6").Value = 5310
$ws.Range("B7").Value = 6.39
$ws.Range("C7").Value = 543
$ws.Range("D7").Value = 8500
$ws.Range("B8").Value = 2.53
$ws.Range("C8").Value = 248
$ws.Range("B9").Value = 1.37
$ws.Range("C9").Value = 134
$ws.Range("B10").Value = 2.15
$ws.Range("C10").Value = 210
